$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textD18 = @'
Aim
id="Par555">To describe associations between frailty, ethnicity, socioeconomic position and mortality in a cohort of older patients presenting to hospital with COVID-19.
Findings
id="Par6666">Frailty did not appear to be associated with mortality rates after COVID-19, though an interaction was evident indicating much larger excess mortality in fitter, compared with frailer patients.


Message
id="Par777">Frailty may not be a good measure of prognosis in COVID-19 and different mechanisms may underlie pathways to death depending on pre-morbid frailty.



'@

$textE18 = @'
[Amy%Miles%NULL%1, Thomas E.%Webb%NULL%2, Thomas E.%Webb%NULL%0, Benjamin C.%Mcloughlin%NULL%2, Benjamin C.%Mcloughlin%NULL%0, Imran%Mannan%NULL%2, Imran%Mannan%NULL%0, Arshad%Rather%NULL%2, Arshad%Rather%NULL%0, Paul%Knopp%NULL%1, Daniel%Davis%daniel.davis@ucl.ac.uk%1]
'@

$textD21 = @'
Background
id="Par1">Understanding of the role of ethnicity and socioeconomic position in the risk of developing SARS-CoV-2 infection is limited.

 We investigated this in the UK Biobank study.


Methods
id="Par2">The UK Biobank study recruited 40–70-year-olds in 2006–2010 from the general population, collecting information about self-defined ethnicity and socioeconomic variables (including area-level socioeconomic deprivation and educational attainment).

 SARS-CoV-2 test results from Public Health England were linked to baseline UK Biobank data.

 Poisson regression with robust standard errors was used to assess risk ratios (RRs) between the exposures and dichotomous variables for being tested, having a positive test and testing positive in hospital.

 We also investigated whether ethnicity and socioeconomic position were associated with having a positive test amongst those tested.

 We adjusted for covariates including age, sex, social variables (including healthcare work and household size), behavioural risk factors and baseline health.


Results
id="Par3">Amongst 392,116 participants in England, 2658 had been tested for SARS-CoV-2 and 948 tested positive (726 in hospital) between 16 March and 3 May 2020. Black and south Asian groups were more likely to test positive (RR 3.35 (95% CI 2.48–4.53) and RR 2.42 (95% CI 1.75–3.36) respectively), with Pakistani ethnicity at highest risk within the south Asian group (RR 3.24 (95% CI 1.73–6.07)).

 These ethnic groups were more likely to be hospital cases compared to the white British.

 Adjustment for baseline health and behavioural risk factors led to little change, with only modest attenuation when accounting for socioeconomic variables.

 Socioeconomic deprivation and having no qualifications were consistently associated with a higher risk of confirmed infection (RR 2.19 for most deprived quartile vs least (95% CI 1.80–2.66) and RR 2.00 for no qualifications vs degree (95% CI 1.66–2.42)).


Conclusions
id="Par4">Some minority ethnic groups have a higher risk of confirmed SARS-CoV-2 infection in the UK Biobank study, which was not accounted for by differences in socioeconomic conditions, baseline self-reported health or behavioural risk factors.

 An urgent response to addressing these elevated risks is required.



'@

$textE21 = @'
[Claire L.%Niedzwiedz%NULL%1, Catherine A.%O’Donnell%NULL%1, Bhautesh Dinesh%Jani%NULL%1, Evangelia%Demou%NULL%1, Frederick K.%Ho%NULL%1, Carlos%Celis-Morales%NULL%1, Barbara I.%Nicholl%NULL%1, Frances S.%Mair%NULL%1, Paul%Welsh%NULL%1, Naveed%Sattar%NULL%1, Jill P.%Pell%NULL%1, S. Vittal%Katikireddi%vittal.katikireddi@glasgow.ac.uk%1]
'@

$textD33 = @'
Background
id="Par1">By mid-May 2020, there were over 1.5 million cases of (SARS-CoV-2) or COVID-19 across the U.

S.

 with new confirmed cases continuing to rise following the re-opening of most states.

 Prior studies have focused mainly on clinical risk factors associated with serious illness and mortality of COVID-19. Less analysis has been conducted on the clinical, sociodemographic, and environmental variables associated with initial infection of COVID-19.
Methods
id="Par2">A multivariable statistical model was used to characterize risk factors in 34,503cases of laboratory-confirmed positive or negative COVID-19 infection in the Providence Health System (U.

S.

) between February 28 and April 27, 2020. Publicly available data were utilized as approximations for social determinants of health, and patient-level clinical and sociodemographic factors were extracted from the electronic medical record.


Results
id="Par3">Higher risk of COVID-19 infection was associated with older age (OR 1.69; 95% CI 1.41–2.02, p &lt; 0.0001), male gender (OR 1.32; 95% CI 1.21–1.44, p &lt; 0.0001), Asian race (OR 1.43; 95% CI 1.18–1.72, p = 0.0002), Black/African American race (OR 1.51; 95% CI 1.25–1.83, p &lt; 0.0001), Latino ethnicity (OR 2.07; 95% CI 1.77–2.41, p &lt; 0.0001), non-English language (OR 2.09; 95% CI 1.7–2.57, p &lt; 0.0001), residing in a neighborhood with financial insecurity (OR 1.10; 95% CI 1.01–1.25, p = 0.04), low air quality (OR 1.01; 95% CI 1.0–1.04, p = 0.05), housing insecurity (OR 1.32; 95% CI 1.16–1.5, p &lt; 0.0001) or transportation insecurity (OR 1.11; 95% CI 1.02–1.23, p = 0.03), and living in senior living communities (OR 1.69; 95% CI 1.23–2.32, p = 0.001).


Conclusion
id="Par4">sisk of COVID-19 infection is higher among groups already affected by health disparities across age, race, ethnicity, language, income, and living conditions.

 Health promotion and disease prevention strategies should prioritize groups most vulnerable to infection and address structural inequities that contribute to risk through social and economic policy.



'@

$textE33 = @'
[Yelena%Rozenfeld%Yelena.Rozenfeld@providence.org%1, Jennifer%Beam%NULL%2, Jennifer%Beam%NULL%0, Haley%Maier%NULL%1, Whitney%Haggerson%NULL%1, Karen%Boudreau%NULL%1, Jamie%Carlson%NULL%1, Rhonda%Medows%NULL%1]
'@

$textD35 = @'
Background
id="Par1">Data from the UK COVID-19 outbreak are emerging, and there are ongoing concerns about a disproportionate effect on ethnic minorities.

 There is very limited information on COVID-19 in the over-80s, and the rates of hospital-onset infections are unknown.


Methods
id="Par2">This was a retrospective cohort study from electronic case records of the first 450 patients admitted to our hospital with PCR-confirmed COVID-19, 77% of the total inpatient caseload to date.

 Demographic, clinical and biochemical data were extracted.

 The primary endpoint was death during the index hospital admission.

 The characteristics of all patients, those over 80 years of age and those with hospital-onset COVID-19 were examined.


Results
id="Par3">The median (IQR) age was 72 (56, 83), with 150 (33%) over 80 years old and 60% male.

 Presenting clinical and biochemical features were consistent with those reported elsewhere.

 The ethnic breakdown of patients admitted was similar to that of our underlying local population.

 Inpatient mortality was high at 38%.


Conclusions
id="Par5">Inpatient mortality was high, especially among the over-80s, who are more likely to present atypically.

 The ethnic composition of our caseload was similar to the underlying population.

 While a significant number of patients acquired COVID-19 while already in hospital, their outcomes were no worse.



'@

$textE35 = @'
[Simon E.%Brill%simon.brill@nhs.net%1, Hannah C.%Jarvis%NULL%2, Hannah C.%Jarvis%NULL%0, Ezgi%Ozcan%NULL%1, Thomas L. P.%Burns%NULL%1, Rabia A.%Warraich%NULL%1, Lisa J.%Amani%NULL%1, Amina%Jaffer%NULL%1, Stephanie%Paget%NULL%1, Anand%Sivaramakrishnan%NULL%1, Dean D.%Creer%NULL%1]
'@

$ws.Cells.Item(18, 4).Value = $textD18
$ws.Cells.Item(18, 5).Value = $textE18
$ws.Cells.Item(21, 4).Value = $textD21
$ws.Cells.Item(21, 5).Value = $textE21
$ws.Cells.Item(33, 4).Value = $textD33
$ws.Cells.Item(33, 5).Value = $textE33
$ws.Cells.Item(35, 4).Value = $textD35
$ws.Cells.Item(35, 5).Value = $textE35
